# Auto-generated Excel COM-interop edit script
# Applies odds corrections to rows 7, 30, 31, 38, 39 and inserts a new match row (new row 46),
# shifting the former row 46 (Tampa Bay vs Hartford Athletic) down to row 47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update individual odds values that changed in existing rows ---
$ws.Range("G7").Value = 3.1
$ws.Range("I7").Value = 2.45
$ws.Range("X7").Value = 13
$ws.Range("AH7").Value = 6.5
$ws.Range("BB7").Value = 301
$ws.Range("G30").Value = 1.27
$ws.Range("O30").Value = 1.14
$ws.Range("P30").Value = 5.5
$ws.Range("S30").Value = 1.25
$ws.Range("T30").Value = 3.75
$ws.Range("U30").Value = 1.91
$ws.Range("V30").Value = 1.91
$ws.Range("W30").Value = 9.5
$ws.Range("X30").Value = 7.5
$ws.Range("Z30").Value = 8.5
$ws.Range("AG30").Value = 251
$ws.Range("AT30").Value = 3.75
$ws.Range("G31").Value = 2.05
$ws.Range("H31").Value = 3.2
$ws.Range("I31").Value = 3.9
$ws.Range("J31").Value = 2.75
$ws.Range("K31").Value = 2.05
$ws.Range("U31").Value = 1.95
$ws.Range("V31").Value = 1.8
$ws.Range("AC31").Value = 8
$ws.Range("AG31").Value = 351
$ws.Range("AH31").Value = 10
$ws.Range("AJ31").Value = 13
$ws.Range("AR31").Value = 67
$ws.Range("AS31").Value = 201
$ws.Range("AZ31").Value = 67
$ws.Range("G38").Value = 2.25
$ws.Range("I38").Value = 2.8
$ws.Range("J38").Value = 2.88
$ws.Range("L38").Value = 3.4
$ws.Range("Q38").Value = 1.67
$ws.Range("R38").Value = 2.15
$ws.Range("S38").Value = 1.3
$ws.Range("T38").Value = 3.4
$ws.Range("X38").Value = 13
$ws.Range("Y38").Value = 9.5
$ws.Range("AB38").Value = 21
$ws.Range("AD38").Value = 7.5
$ws.Range("AE38").Value = 12
$ws.Range("AG38").Value = 126
$ws.Range("AK38").Value = 29
$ws.Range("AT38").Value = 3.4
$ws.Range("AU38").Value = 7
$ws.Range("BA38").Value = 51
$ws.Range("G39").Value = 2.8
$ws.Range("I39").Value = 2.63
$ws.Range("L39").Value = 3.4
$ws.Range("Y39").Value = 11
$ws.Range("AM39").Value = 41
$ws.Range("AY39").Value = 29
$ws.Range("BB39").Value = 251

# --- Step 2: insert a new row at position 46, pushing the existing row 46 down to row 47 ---
$ws.Rows.Item(46).Insert()

# --- Step 3: populate the newly inserted row 46 with the new match data ---
$ws.Range("A46").Value = 'jJhlHOOk'
$ws.Range("B46").Value = '23/10/2024'
$ws.Range("C46").Value = '23:30'
$ws.Range("D46").Value = 'USA - MLS'
$ws.Range("E46").Value = 'Vancouver Whitecaps'
$ws.Range("F46").Value = 'Portland Timbers'
$ws.Range("G46").Value = 3.4
$ws.Range("H46").Value = 4
$ws.Range("I46").Value = 1.95
$ws.Range("J46").Value = 3.75
$ws.Range("K46").Value = 2.4
$ws.Range("L46").Value = 2.5
$ws.Range("M46").Value = 1.02
$ws.Range("N46").Value = 17
$ws.Range("O46").Value = 1.13
$ws.Range("P46").Value = 5
$ws.Range("Q46").Value = 1.53
$ws.Range("R46").Value = 2.4
$ws.Range("S46").Value = 1.25
$ws.Range("T46").Value = 3.75
$ws.Range("U46").Value = 1.47
$ws.Range("V46").Value = 2.5
$ws.Range("W46").Value = 15
$ws.Range("X46").Value = 21
$ws.Range("Y46").Value = 12
$ws.Range("Z46").Value = 41
$ws.Range("AA46").Value = 23
$ws.Range("AB46").Value = 26
$ws.Range("AC46").Value = 19
$ws.Range("AD46").Value = 8
$ws.Range("AE46").Value = 12
$ws.Range("AF46").Value = 34
$ws.Range("AG46").Value = 101
$ws.Range("AH46").Value = 11
$ws.Range("AI46").Value = 12
$ws.Range("AJ46").Value = 9
$ws.Range("AK46").Value = 19
$ws.Range("AL46").Value = 13
$ws.Range("AM46").Value = 19
$ws.Range("AN46").Value = 6
$ws.Range("AO46").Value = 17
$ws.Range("AP46").Value = 21
$ws.Range("AQ46").Value = 51
$ws.Range("AR46").Value = 51
$ws.Range("AS46").Value = 101
$ws.Range("AT46").Value = 3.75
$ws.Range("AU46").Value = 7
$ws.Range("AV46").Value = 41
$ws.Range("AW46").Value = 4.33
$ws.Range("AX46").Value = 10
$ws.Range("AY46").Value = 17
$ws.Range("AZ46").Value = 29
$ws.Range("BA46").Value = 41
$ws.Range("BB46").Value = 81
$ws.Range("BC46").Value = 301
$ws.Range("BD46").Value = 151
